# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" right before "总计", formatted like the
#    most recent quarter sheet ("2021-Q4"), populated with the new quarter's
#    per-fund holdings.
# 2. Prepend a new "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing rows down and renumbering the helper index column.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1. New "2022-Q1" sheet, inserted before "总计" -----------------------
# NOTE: Worksheets.Add(before) re-aliases the handle passed as the "before"
# argument to the newly created sheet, so grab a *fresh* handle to "总计" by
# name afterwards rather than reusing the one passed into Add().
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

# Match the page margins used throughout the rest of this workbook (values
# are in points: 0.75in/1in/0.5in).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Copy header row + row-2 formatting/layout from 2021-Q4 so styles (bold
# centered header, bordered index column, etc.) match the existing sheets.
$q4.Range("B1:H1").Copy($q1.Range("B1"))
$q4.Range("A2:H2").Copy($q1.Range("A2"))
$q4.Range("A2:H2").Copy($q1.Range("A3"))
$q4.Range("A2:H2").Copy($q1.Range("A4"))

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2: 富国生物医药科技混合A
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "006218"
Set-TextValue $q1.Range("C2") "富国生物医药科技混合A"
Set-TextValue $q1.Range("D2") "9.55"
Set-TextValue $q1.Range("E2") "87.25"
Set-TextValue $q1.Range("F2") "5.33"
Set-TextValue $q1.Range("G2") "0.5090"
$q1.Range("H2").Value = 3

# Row 3: 富国天源沪港深平衡混合
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "100016"
Set-TextValue $q1.Range("C3") "富国天源沪港深平衡混合"
Set-TextValue $q1.Range("D3") "6.23"
Set-TextValue $q1.Range("E3") "72.29"
Set-TextValue $q1.Range("F3") "3.01"
Set-TextValue $q1.Range("G3") "0.1875"
$q1.Range("H3").Value = 6

# Row 4: 富国生物医药科技混合C
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "011308"
Set-TextValue $q1.Range("C4") "富国生物医药科技混合C"
Set-TextValue $q1.Range("D4") "1.26"
Set-TextValue $q1.Range("E4") "87.25"
Set-TextValue $q1.Range("F4") "5.33"
Set-TextValue $q1.Range("G4") "0.0672"
$q1.Range("H4").Value = 3

# --- 2. Prepend a "2022-Q1" row to "总计" ----------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()
$zj.Range("A3:D3").Copy($zj.Range("A2"))

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 3
$zj.Range("D2").Value = 0.76

$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3
$zj.Range("A6").Value = 4

# Keep the originally active tab ("2021-Q1") selected, same as before the edit.
$wb.Worksheets.Item("2021-Q1").Select()
